{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// This document is a \"convention de stage\" generated from a template. The\n// commit only changes two pieces of visible text:\n//   1. \"Entre  2025\"                       -> \"Entre  2024\"\n//   2. \"FAIT \u00e0 MONTREUIL LE : 11/01/2025\"  -> \"FAIT \u00e0 MONTREUIL LE : 26/02/2025\"\n// (the surrounding \"Dates : Du 2025-01-15 Au 2025-07-15\" text elsewhere in\n// the document is untouched by the diff, so we must match precisely rather\n// than doing a blind \"2025\" -> \"2024\" replace).\n\n// 1) \"Entre  2025\" -> \"Entre  2024\" (two spaces, exact phrase is unique in doc).\nconst introResults = context.document.body.search(\"Entre  2025\", { matchCase: true, matchWholeWord: false });\nintroResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < introResults.items.length; i++) {\n  introResults.items[i].insertText(\"Entre  2024\", Word.InsertLocation.replace);\n}\n\n// 2) Update the signature date. Replace only the date token so we don't have\n// to worry about the exact (non-breaking) space before the colon.\nconst dateResults = context.document.body.search(\"11/01/2025\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"26/02/2025\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# This document is a \"convention de stage\" generated from a template. The\n# commit only changes two pieces of visible text:\n#   1. \"Entre  2025\"                       -> \"Entre  2024\"\n#   2. \"FAIT \u00e0 MONTREUIL LE : 11/01/2025\"  -> \"FAIT \u00e0 MONTREUIL LE : 26/02/2025\"\n# (the unrelated \"Dates : Du 2025-01-15 Au 2025-07-15\" text elsewhere in the\n# document must stay untouched, so we match precise, unique phrases instead\n# of a blind \"2025\" -> \"2024\" replace).\n\n$d = $word.ActiveDocument\n\n# 1) \"Entre  2025\" -> \"Entre  2024\" (two spaces; unique in the document).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Entre  2025\", $true, $false, $false, $false, $false, $true, 1, $false, \"Entre  2024\", 2) | Out-Null\n\n# 2) Update the signature date. Replace only the date token so the\n# (non-breaking) space before the colon is left exactly as it was.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"11/01/2025\", $true, $false, $false, $false, $false, $true, 1, $false, \"26/02/2025\", 2) | Out-Null\n\n$d.Saved = $false\n"}
